# Commit: "Fruta / hortaliza, semanal"
#
# A brand-new weekly price report row is inserted into the data table at
# row 49 (just after the existing "Provincia de Quillota" bandeja record),
# pushing every following record down by one row (old row 49 -> new row 50,
# old row 121 -> new row 122). The sheet's dimension grows from A1:T121 to
# A1:T122 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 49; Excel shifts rows 49:121 down to 50:122
# and carries the existing row formatting (e.g. the date style on column D)
# along with it.
$ws.Rows("49:49").Insert()

# Populate the newly inserted row with this week's report.
$ws.Range("A49").Value = 10
$ws.Range("B49").Value = "Vega Modelo de Temuco"
$ws.Range("C49").Value = "La Araucanía"
$ws.Range("D49").Value = 44915
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100101
$ws.Range("H49").Value = "Berries"
$ws.Range("I49").Value = 100101001
$ws.Range("J49").Value = "Arándano (blue)"
$ws.Range("K49").Value = "Sin especificar"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 500
$ws.Range("N49").Value = 2000
$ws.Range("O49").Value = 2000
$ws.Range("P49").Value = 2000
$ws.Range("Q49").Value = "$/kilo"
$ws.Range("R49").Value = "Región del Maule"
$ws.Range("S49").Value = 2000
$ws.Range("T49").Value = 1
